$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Duplicate rows 2-10 (columns C:E) into new rows 11-19 on sheet1
# ("added try/except for openpyxl-max_row issue" -> duplicate the data block
# so downstream tooling that relied on max_row still finds rows after the fix)
for ($i = 0; $i -le 8; $i++) {
    $srcRow = 2 + $i
    $dstRow = 11 + $i

    $srcRange = $ws1.Range("C" + $srcRow + ":E" + $srcRow)
    $dstRange = $ws1.Range("C" + $dstRow + ":E" + $dstRow)

    # copy values first (keeps shared-string reuse), then formats
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4163)

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Update the selections / active sheet to match the saved view state
$ws1.Activate()
$ws1.Range("C11:E19").Select()

$ws2.Activate()
$ws2.Range("B13").Select()

$ws1.Activate()
